$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 27 de Junio de 2020 a las 12:23"

# Row 7: India
$ws.Range("A7").Value = "India"
$ws.Range("B7").Value = 510672
$ws.Range("C7").Value = 1226
$ws.Range("D7").Value = 296376
$ws.Range("E7").Value = 198584
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 23
$ws.Range("H7").Value = 15712

# Row 28: Belgica
$ws.Range("A28").Value = "Belgica"
$ws.Range("B28").Value = 61209
$ws.Range("C28").Value = 103
$ws.Range("D28").Value = 16941
$ws.Range("E28").Value = 34536
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 9732

# Row 40: Oman
$ws.Range("A40").Value = "Oman"
$ws.Range("B40").Value = 36953
$ws.Range("C40").Value = 919
$ws.Range("D40").Value = 20363
$ws.Range("E40").Value = 16431
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 6
$ws.Range("H40").Value = 159

# Row 67: Marruecos
$ws.Range("A67").Value = "Marruecos"
$ws.Range("B67").Value = 11854
$ws.Range("C67").Value = 221
$ws.Range("D67").Value = 8700
$ws.Range("E67").Value = 2936
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 218

# Row 68: Nepal
$ws.Range("A68").Value = "Nepal"
$ws.Range("B68").Value = 11755
$ws.Range("C68").Value = 0
$ws.Range("D68").Value = 2698
$ws.Range("E68").Value = 9030
$ws.Range("F68").Value = 0
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 27

# Row 73: Malasia
$ws.Range("A73").Value = "Malasia"
$ws.Range("B73").Value = 8616
$ws.Range("C73").Value = 10
$ws.Range("D73").Value = 8308
$ws.Range("E73").Value = 187
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = 121

# Row 76: Finlandia
$ws.Range("A76").Value = "Finlandia"
$ws.Range("B76").Value = 7198
$ws.Range("C76").Value = 7
$ws.Range("D76").Value = 6600
$ws.Range("E76").Value = 270
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 0
$ws.Range("H76").Value = 328

# Row 103: Albania
$ws.Range("A103").Value = "Albania"
$ws.Range("B103").Value = 2330
$ws.Range("C103").Value = 61
$ws.Range("D103").Value = 1346
$ws.Range("E103").Value = 931
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 2
$ws.Range("H103").Value = 53

# Row 104: Cuba
$ws.Range("A104").Value = "Cuba"
$ws.Range("B104").Value = 2325
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 2180
$ws.Range("E104").Value = 60
$ws.Range("F104").Value = 0
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 85

# Row 105: Maldivas
$ws.Range("A105").Value = "Maldivas"
$ws.Range("B105").Value = 2283
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 1863
$ws.Range("E105").Value = 412
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 8

# Row 108: Sri Lanka
$ws.Range("A108").Value = "Sri Lanka"
$ws.Range("B108").Value = 2033
$ws.Range("C108").Value = 19
$ws.Range("D108").Value = 1639
$ws.Range("E108").Value = 383
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 11

# Row 124: Hong Kong
$ws.Range("A124").Value = "Hong Kong"
$ws.Range("B124").Value = 1198
$ws.Range("C124").Value = 1
$ws.Range("D124").Value = 1095
$ws.Range("E124").Value = 96
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 7

# Row 171: Gibraltar
$ws.Range("A171").Value = "Gibraltar"
$ws.Range("B171").Value = 177
$ws.Range("C171").Value = 1
$ws.Range("D171").Value = 176
$ws.Range("E171").Value = 1
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 0

# Row 201: Santa Lucia
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("B201").Value = 19
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 19
$ws.Range("E201").Value = 0
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0

# Row 202: Laos
$ws.Range("A202").Value = "Laos"
$ws.Range("B202").Value = 19
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 19
$ws.Range("E202").Value = 0
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

# Row 203: Fiyi
$ws.Range("A203").Value = "Fiyi"
$ws.Range("B203").Value = 18
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 18
$ws.Range("E203").Value = 0
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

# Row 204: Dominica
$ws.Range("A204").Value = "Dominica"
$ws.Range("B204").Value = 18
$ws.Range("C204").Value = 0
$ws.Range("D204").Value = 18
$ws.Range("E204").Value = 0
$ws.Range("F204").Value = 0
$ws.Range("G204").Value = 0
$ws.Range("H204").Value = 0

# Row 208: Groenlandia
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("B208").Value = 13
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 13
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0

# Row 209: Islas Malvinas
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("B209").Value = 13
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 13
$ws.Range("E209").Value = 0
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

# Row 212: Seychelles
$ws.Range("A212").Value = "Seychelles"
$ws.Range("B212").Value = 11
$ws.Range("C212").Value = 0
$ws.Range("D212").Value = 11
$ws.Range("E212").Value = 0
$ws.Range("F212").Value = 0
$ws.Range("G212").Value = 0
$ws.Range("H212").Value = 0

# Row 213: Montserrat
$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 11
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 10
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

